# Regenerate the "K" column (G) values for the save_data sheet.
# The commit switches the K column computation away from the old
# "Strike#" based value to the newly calculated s_vals, updating the
# stored constants for rows 2-11 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 3
    6  = 0
    7  = 0
    8  = 1
    9  = 2
    10 = 3
    11 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
